# Generate Report for Handoff
# Regenerates the localization-status report with a fresh handoff:
#   - old source files 61df56d7...md / b3c9e954...md
#   - replaced by 2246d8cb...md / ffff7e577d15...md
#   - status moves from "Handed back: in sync with en-US" -> "Ready for handoff"
#   - the old duplicate "legacy" F/G columns on the per-locale sheets are dropped

function Remove-HyperlinkAt($ws, $addr) {
    $found = $null
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $addr) {
            $found = $h
            break
        }
    }
    if ($found -ne $null) {
        $found.Delete()
    }
}

function Set-HyperlinkDisplay($ws, $addr, $text) {
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $addr) {
            $h.TextToDisplay = $text
            break
        }
    }
}

$wb = $excel.ActiveWorkbook

$file1Md   = "2246d8cb-028b-463f-8a0a-0d8d45762021.md"
$file2Md   = "ffff7e577d15-3540-4630-af39-f5b803a4b64e.md"
$status    = "Ready for handoff"
$overviewDate = "2016-03-13 11:03:45"

$zhXlf = "2246d8cb-028b-463f-8a0a-0d8d45762021.1c8bef3171a484bb26a51d73e8ed7926ea9ada6e.zh-cn.xlf"
$deXlf = "2246d8cb-028b-463f-8a0a-0d8d45762021.1c8bef3171a484bb26a51d73e8ed7926ea9ada6e.de-de.xlf"
$zhHandoffDatetime = "2016-03-13 11:03:41"
$deHandoffDatetime = "2016-03-13 11:03:45"
$neverHandback = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = $file1Md
$ws1.Range("B2").Value = $status
$ws1.Range("C2").Value = $status
$ws1.Range("D2").Value = $overviewDate

$ws1.Range("A3").Value = $file2Md
$ws1.Range("B3").Value = $status
$ws1.Range("C3").Value = $status
$ws1.Range("D3").Value = $overviewDate

Set-HyperlinkDisplay $ws1 '$A$2' $file1Md
Set-HyperlinkDisplay $ws1 '$A$3' $file2Md

# ---------------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = $file1Md
$ws2.Range("C2").Value = $status
$ws2.Range("D2").Value = $zhXlf
$ws2.Range("E2").Value = $zhHandoffDatetime
$ws2.Range("H2").Value = $neverHandback

$ws2.Range("A3").Value = $file2Md
$ws2.Range("C3").Value = $status
$ws2.Range("D3").Value = $zhXlf
$ws2.Range("E3").Value = $zhHandoffDatetime
$ws2.Range("H3").Value = $neverHandback

Remove-HyperlinkAt $ws2 '$F$2'
Remove-HyperlinkAt $ws2 '$G$2'
Remove-HyperlinkAt $ws2 '$F$3'
Remove-HyperlinkAt $ws2 '$G$3'
$ws2.Range("F2:G3").Clear()

Set-HyperlinkDisplay $ws2 '$A$2' $file1Md
Set-HyperlinkDisplay $ws2 '$D$2' $zhXlf
Set-HyperlinkDisplay $ws2 '$A$3' $file2Md
Set-HyperlinkDisplay $ws2 '$D$3' $zhXlf

# ---------------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = $file1Md
$ws3.Range("C2").Value = $status
$ws3.Range("D2").Value = $deXlf
$ws3.Range("E2").Value = $deHandoffDatetime
$ws3.Range("H2").Value = $neverHandback

$ws3.Range("A3").Value = $file2Md
$ws3.Range("C3").Value = $status
$ws3.Range("D3").Value = $deXlf
$ws3.Range("E3").Value = $deHandoffDatetime
$ws3.Range("H3").Value = $neverHandback

Remove-HyperlinkAt $ws3 '$F$2'
Remove-HyperlinkAt $ws3 '$G$2'
Remove-HyperlinkAt $ws3 '$F$3'
Remove-HyperlinkAt $ws3 '$G$3'
$ws3.Range("F2:G3").Clear()

Set-HyperlinkDisplay $ws3 '$A$2' $file1Md
Set-HyperlinkDisplay $ws3 '$D$2' $deXlf
Set-HyperlinkDisplay $ws3 '$A$3' $file2Md
Set-HyperlinkDisplay $ws3 '$D$3' $deXlf
